$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.239.44'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '2.305.56'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("E4").Value = '  -0.04%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '302.16'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '100.15'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +2.89%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.506'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  +3.71%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '36.45'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +7.96%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0794'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.02%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '18.71'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +10.40%  '

$ws.Range("E13").Value = '  +1.02%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '7.01'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +3.17%  '

$ws.Range("D15").Value = '2.663.45'
$ws.Range("E15").Value = '  -0.04%  '

$ws.Range("D16").Value = '2.351.88'
$ws.Range("E16").Value = '  +3.19%  '

$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("D18").Value = '43.115.90'
$ws.Range("E18").Value = '  +0.27%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '12.84'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +10.26%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '6.19'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +2.88%  '

$ws.Range("D21").Value = '0.0₃0908'
$ws.Range("E21").Value = '  +0.70%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '68.07'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +1.30%  '

$ws.Range("E23").Value = '  +13.90%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '236.47'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.33%  '

$ws.Range("E26").Value = '  -0.40%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '25.22'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.35%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.35'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +9.04%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '34.84'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +2.10%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '167.26'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.29%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '9.18'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.41%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.05%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.06'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +1.74%  '

$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '17.86'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +5.56%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '4.72'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -2.04%  '

$ws.Range("E36").Value = '  +0.81%  '

$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("E38").Value = '  +0.30%  '

$ws.Range("E39").Value = '  +2.30%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.102'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +0.75%  '

$ws.Range("E41").Value = '  +0.73%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '2.34'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("D43").Value = '1.990.33'
$ws.Range("E43").Value = '  +0.82%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0292'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +3.62%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '10.13'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +3.42%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '17.67'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.12%  '

$ws.Range("E47").Value = '  +1.56%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '55.65'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +5.14%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '1.55'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +3.37%  '

$ws.Range("D50").Value = '2.531.31'
$ws.Range("E50").Value = '  +0.12%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '71.02'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +1.50%  '

